# FAST_holdings.xlsx - "Add files via upload"
#
# The model-holdings date stamp moved from 2021-07-08 to 2021-07-09 and the
# Weight (D) / Percent Change (E) figures for rows 2-10 were refreshed to
# match. The sheet ships protected (no UI password known to us), so it has
# to be unprotected before the cells can be written and is re-protected
# afterwards to leave the workbook in the same guarded state it started in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Unprotect()

# Update the confidentiality / "as of" date footnote (shared string used by A13)
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-09 for illustrative purposes only and are subject to change."

# Row 2
$ws.Range("D2").Value = 0.1369940756157144
$ws.Range("E2").Value = 0.01477221470455548

# Row 3
$ws.Range("D3").Value = 0.1078948825013534
$ws.Range("E3").Value = 0.02054673200684887

# Row 4
$ws.Range("D4").Value = 0.1131818081883071
$ws.Range("E4").Value = 0.0088000596614215

# Row 5
$ws.Range("D5").Value = 0.1187604966206444
$ws.Range("E5").Value = 0.0148946099401106

# Row 6
$ws.Range("D6").Value = 0.1221582903830675
$ws.Range("E6").Value = 0.006982677588290631

# Row 7
$ws.Range("D7").Value = 0.1424627159819709
$ws.Range("E7").Value = 0.01731180704143154

# Row 8
$ws.Range("D8").Value = 0.1302325695469338
$ws.Range("E8").Value = 0.01525423728813546

# Row 9
$ws.Range("D9").Value = 0.1283151611620085
$ws.Range("E9").Value = 0.0142170951995868

# Row 10 (Total)
$ws.Range("D10").Value = 0.9999999999999999
$ws.Range("E10").Value = 0.01413563745847912

# Restore the sheet protection the workbook shipped with.
$ws.Protect([System.Reflection.Missing]::Value, $true, $true, $true)
